# Insert a new data row at row 440, shifting the existing rows 440:526 down
# to 441:527, and populate the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 440 (shifts rows 440:526 -> 441:527)
$ws.Rows("440:440").Insert()

# Populate the new row 440 with the same "fixed" columns used throughout the
# table plus the specific values for this record.
$ws.Range("A440").Value = 10
$ws.Range("B440").Value = "Vega Modelo de Temuco"
$ws.Range("C440").Value = "La Araucanía"
$ws.Range("D440").Value = (Get-Date -Year 2023 -Month 3 -Day 10).Date
$ws.Range("E440").Value = 9
$ws.Range("F440").Value = "Fruta"
$ws.Range("G440").Value = 100108
$ws.Range("H440").Value = "Tropicales y subtropicales"
$ws.Range("I440").Value = 100108002
$ws.Range("J440").Value = "Mango"
$ws.Range("K440").Value = "Sin especificar"
$ws.Range("L440").Value = "Primera"
$ws.Range("M440").Value = 500
$ws.Range("N440").Value = 8000
$ws.Range("O440").Value = 8000
$ws.Range("P440").Value = 8000
$ws.Range("Q440").Value = "`$/bandeja 4 kilos"
$ws.Range("R440").Value = "Perú"
$ws.Range("S440").Value = 2000
$ws.Range("T440").Value = 4
